# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# F2: 1417 -> 1418
# F3: 2980 -> 2981
# These values live on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1418
    $ws.Range("F3").Value = 2981
}
